$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Template rows for formatting: row 2 has column N populated (L=Black), row 3 has column M populated (L=Red)
$templateRowN = 2
$templateRowM = 3

# Row 617
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A617:N617").PasteSpecial(-4122) | Out-Null
$ws.Range("A617").Value = 45191.93974084491
$ws.Range("B617").Value = 'godis2002@naver.com'
$ws.Range("C617").Value = '러시아학과'
$ws.Range("D617").Value = 20211706
$ws.Range("E617").Value = '김승겸'
$ws.Range("F617").Value = '74:26'
$ws.Range("G617").Value = 0.2
$ws.Range("H617").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I617").Value = '952만 명'
$ws.Range("J617").Value = 0.059
$ws.Range("K617").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L617").Value = 'Red'
$ws.Range("M617").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F617").PasteSpecial(-4122) | Out-Null

# Row 618
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A618:N618").PasteSpecial(-4122) | Out-Null
$ws.Range("A618").Value = 45191.9549290625
$ws.Range("B618").Value = 'juh0611@naver.com'
$ws.Range("C618").Value = '경영학과'
$ws.Range("D618").Value = 20221051
$ws.Range("E618").Value = '신주희'
$ws.Range("F618").Value = '74:26'
$ws.Range("G618").Value = 0.2
$ws.Range("H618").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I618").Value = '952만 명'
$ws.Range("J618").Value = 0.059
$ws.Range("K618").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L618").Value = 'Black'
$ws.Range("N618").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F618").PasteSpecial(-4122) | Out-Null

# Row 619
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A619:N619").PasteSpecial(-4122) | Out-Null
$ws.Range("A619").Value = 45191.973035
$ws.Range("B619").Value = 'rlagkdud114@naver.com'
$ws.Range("C619").Value = '정치행정학과'
$ws.Range("D619").Value = 20221026
$ws.Range("E619").Value = '김하영'
$ws.Range("F619").Value = '77:23'
$ws.Range("G619").Value = 0.2
$ws.Range("H619").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I619").Value = '779만 명'
$ws.Range("J619").Value = 0.151
$ws.Range("K619").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L619").Value = 'Black'
$ws.Range("N619").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F619").PasteSpecial(-4122) | Out-Null

# Row 620
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A620:N620").PasteSpecial(-4122) | Out-Null
$ws.Range("A620").Value = 45191.98110314815
$ws.Range("B620").Value = 'sungjuwon1@gmail.com'
$ws.Range("C620").Value = '컨텐츠IT전공'
$ws.Range("D620").Value = 20215173
$ws.Range("E620").Value = '성주원'
$ws.Range("F620").Value = '74:26'
$ws.Range("G620").Value = 0.2
$ws.Range("H620").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I620").Value = '952만 명'
$ws.Range("J620").Value = 0.059
$ws.Range("K620").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L620").Value = 'Black'
$ws.Range("N620").Value = '모름/무응답'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F620").PasteSpecial(-4122) | Out-Null

# Row 621
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A621:N621").PasteSpecial(-4122) | Out-Null
$ws.Range("A621").Value = 45191.98191899306
$ws.Range("B621").Value = 'janghangyeol0304@gmail.com'
$ws.Range("C621").Value = '경제학과'
$ws.Range("D621").Value = 20222838
$ws.Range("E621").Value = '장한결'
$ws.Range("F621").Value = '74:26'
$ws.Range("G621").Value = 0.2
$ws.Range("H621").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I621").Value = '952만 명'
$ws.Range("J621").Value = 0.374
$ws.Range("K621").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L621").Value = 'Red'
$ws.Range("M621").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F621").PasteSpecial(-4122) | Out-Null

# Row 622
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A622:N622").PasteSpecial(-4122) | Out-Null
$ws.Range("A622").Value = 45192.003076678244
$ws.Range("B622").Value = '20233814@hallym.ac.kr'
$ws.Range("C622").Value = '식품영양학과'
$ws.Range("D622").Value = 20233814
$ws.Range("E622").Value = '김정현'
$ws.Range("F622").Value = '74:26'
$ws.Range("G622").Value = 0.2
$ws.Range("H622").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I622").Value = '952만 명'
$ws.Range("J622").Value = 0.059
$ws.Range("K622").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L622").Value = 'Black'
$ws.Range("N622").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F622").PasteSpecial(-4122) | Out-Null

# Row 623
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A623:N623").PasteSpecial(-4122) | Out-Null
$ws.Range("A623").Value = 45192.13073451389
$ws.Range("B623").Value = 'ridsigdog@gmail.com'
$ws.Range("C623").Value = '경제학과'
$ws.Range("D623").Value = 20202816
$ws.Range("E623").Value = '박건민'
$ws.Range("F623").Value = '74:26'
$ws.Range("G623").Value = 0.25
$ws.Range("H623").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I623").Value = '779만 명'
$ws.Range("J623").Value = 0.151
$ws.Range("K623").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("L623").Value = 'Black'
$ws.Range("N623").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F623").PasteSpecial(-4122) | Out-Null

# Row 624
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A624:N624").PasteSpecial(-4122) | Out-Null
$ws.Range("A624").Value = 45192.39833813657
$ws.Range("B624").Value = 'jimin4729@naver.com'
$ws.Range("C624").Value = '법학과'
$ws.Range("D624").Value = 20232747
$ws.Range("E624").Value = '임지민'
$ws.Range("F624").Value = '75:25'
$ws.Range("G624").Value = 0.15
$ws.Range("H624").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I624").Value = '166만 명'
$ws.Range("J624").Value = 0.374
$ws.Range("K624").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L624").Value = 'Red'
$ws.Range("M624").Value = '모름/무응답'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F624").PasteSpecial(-4122) | Out-Null

# Row 625
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A625:N625").PasteSpecial(-4122) | Out-Null
$ws.Range("A625").Value = 45192.49826432871
$ws.Range("B625").Value = 'jsy5233406@naver.com'
$ws.Range("C625").Value = '식품영양학과'
$ws.Range("D625").Value = 20233844
$ws.Range("E625").Value = '정서연'
$ws.Range("F625").Value = '75:25'
$ws.Range("G625").Value = 0.15
$ws.Range("H625").Value = '프랑스와 스웨덴의 국민부담률은 꾸준히 40%를 넘고 있다.'
$ws.Range("I625").Value = '952만 명'
$ws.Range("J625").Value = 0.059
$ws.Range("K625").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L625").Value = 'Red'
$ws.Range("M625").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F625").PasteSpecial(-4122) | Out-Null

# Row 626
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A626:N626").PasteSpecial(-4122) | Out-Null
$ws.Range("A626").Value = 45192.5085250463
$ws.Range("B626").Value = '040415kimdh@naver.com'
$ws.Range("C626").Value = '소프트웨어학부'
$ws.Range("D626").Value = 20235119
$ws.Range("E626").Value = '김대현'
$ws.Range("F626").Value = '74:26'
$ws.Range("G626").Value = 0.2
$ws.Range("H626").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I626").Value = '952만 명'
$ws.Range("J626").Value = 0.059
$ws.Range("K626").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L626").Value = 'Black'
$ws.Range("N626").Value = '모름/무응답'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F626").PasteSpecial(-4122) | Out-Null

# Row 627
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A627:N627").PasteSpecial(-4122) | Out-Null
$ws.Range("A627").Value = 45192.53221329861
$ws.Range("B627").Value = '20226760@hallym.ac.kr'
$ws.Range("C627").Value = 'Ai 의료융합'
$ws.Range("D627").Value = 20226760
$ws.Range("E627").Value = '이민홍'
$ws.Range("F627").Value = '74:26'
$ws.Range("G627").Value = 0.2
$ws.Range("H627").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I627").Value = '952만 명'
$ws.Range("J627").Value = 0.059
$ws.Range("K627").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L627").Value = 'Red'
$ws.Range("M627").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F627").PasteSpecial(-4122) | Out-Null

# Row 628
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A628:N628").PasteSpecial(-4122) | Out-Null
$ws.Range("A628").Value = 45192.54229869213
$ws.Range("B628").Value = 'ans1929@gmail.com'
$ws.Range("C628").Value = '법학과'
$ws.Range("D628").Value = 20202750
$ws.Range("E628").Value = '조유진'
$ws.Range("F628").Value = '74:26'
$ws.Range("G628").Value = 0.2
$ws.Range("H628").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I628").Value = '952만 명'
$ws.Range("J628").Value = 0.059
$ws.Range("K628").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L628").Value = 'Red'
$ws.Range("M628").Value = '모름/무응답'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F628").PasteSpecial(-4122) | Out-Null

# Row 629
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A629:N629").PasteSpecial(-4122) | Out-Null
$ws.Range("A629").Value = 45192.54367741898
$ws.Range("B629").Value = 'a01035025756@gmail.com'
$ws.Range("C629").Value = '사회복지학부'
$ws.Range("D629").Value = 20232342
$ws.Range("E629").Value = '이건희'
$ws.Range("F629").Value = '74:26'
$ws.Range("G629").Value = 0.2
$ws.Range("H629").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I629").Value = '779만 명'
$ws.Range("J629").Value = 0.059
$ws.Range("K629").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L629").Value = 'Black'
$ws.Range("N629").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F629").PasteSpecial(-4122) | Out-Null

# Row 630
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A630:N630").PasteSpecial(-4122) | Out-Null
$ws.Range("A630").Value = 45192.56995640046
$ws.Range("B630").Value = 'kimhongik03@naver.com'
$ws.Range("C630").Value = '심리학과'
$ws.Range("D630").Value = 20232114
$ws.Range("E630").Value = '김홍익'
$ws.Range("F630").Value = '76:24'
$ws.Range("G630").Value = 0.2
$ws.Range("H630").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I630").Value = '779만 명'
$ws.Range("J630").Value = 0.374
$ws.Range("K630").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L630").Value = 'Red'
$ws.Range("M630").Value = '모름/무응답'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F630").PasteSpecial(-4122) | Out-Null

# Row 631
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A631:N631").PasteSpecial(-4122) | Out-Null
$ws.Range("A631").Value = 45192.58355498842
$ws.Range("B631").Value = 'hyunjong9951@gmail.com'
$ws.Range("C631").Value = '디스플레이'
$ws.Range("D631").Value = 20183319
$ws.Range("E631").Value = '이현종'
$ws.Range("F631").Value = '74:26'
$ws.Range("G631").Value = 0.2
$ws.Range("H631").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I631").Value = '952만 명'
$ws.Range("J631").Value = 0.059
$ws.Range("K631").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L631").Value = 'Black'
$ws.Range("N631").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F631").PasteSpecial(-4122) | Out-Null

# Row 632
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A632:N632").PasteSpecial(-4122) | Out-Null
$ws.Range("A632").Value = 45192.59258246528
$ws.Range("B632").Value = 'serf0403@naver.com'
$ws.Range("C632").Value = '바이오메디컬학과'
$ws.Range("D632").Value = 20193646
$ws.Range("E632").Value = '정예선'
$ws.Range("F632").Value = '77:23'
$ws.Range("G632").Value = 0.2
$ws.Range("H632").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I632").Value = '166만 명'
$ws.Range("J632").Value = 0.151
$ws.Range("K632").Value = '상호출자제한기업은 신고법인수의 0.1%를 차지하는 데 부담하는 세액은 25.5%이다'
$ws.Range("L632").Value = 'Red'
$ws.Range("M632").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F632").PasteSpecial(-4122) | Out-Null

# Row 633
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A633:N633").PasteSpecial(-4122) | Out-Null
$ws.Range("A633").Value = 45192.598411886574
$ws.Range("B633").Value = 'wnruddms@naver.com'
$ws.Range("C633").Value = '식품영양학과'
$ws.Range("D633").Value = 20203842
$ws.Range("E633").Value = '주경은'
$ws.Range("F633").Value = '74:26'
$ws.Range("G633").Value = 0.2
$ws.Range("H633").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I633").Value = '952만 명'
$ws.Range("J633").Value = 0.059
$ws.Range("K633").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L633").Value = 'Black'
$ws.Range("N633").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F633").PasteSpecial(-4122) | Out-Null

# Row 634
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A634:N634").PasteSpecial(-4122) | Out-Null
$ws.Range("A634").Value = 45192.61100075232
$ws.Range("B634").Value = 'freenix1001@naver.com'
$ws.Range("C634").Value = '언어청각학부'
$ws.Range("D634").Value = 20233903
$ws.Range("E634").Value = '강서연'
$ws.Range("F634").Value = '74:26'
$ws.Range("G634").Value = 0.2
$ws.Range("H634").Value = 'OECD평균은 2010년 31.6%에서 2020년 33.5%까지 상승하였다.'
$ws.Range("I634").Value = '779만 명'
$ws.Range("J634").Value = 0.002
$ws.Range("K634").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L634").Value = 'Red'
$ws.Range("M634").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F634").PasteSpecial(-4122) | Out-Null

# Row 635
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A635:N635").PasteSpecial(-4122) | Out-Null
$ws.Range("A635").Value = 45192.61873761574
$ws.Range("B635").Value = 'jonggwang0104@naver.com'
$ws.Range("C635").Value = '식품영양학과'
$ws.Range("D635").Value = 20183820
$ws.Range("E635").Value = '박종광'
$ws.Range("F635").Value = '74:26'
$ws.Range("G635").Value = 0.2
$ws.Range("H635").Value = '미국의 국민부담률은 우리나라의 국민부담률보다 항상 높은 수준을 유지하고 있다.'
$ws.Range("I635").Value = '952만 명'
$ws.Range("J635").Value = 0.059
$ws.Range("K635").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L635").Value = 'Red'
$ws.Range("M635").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F635").PasteSpecial(-4122) | Out-Null

# Row 636
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A636:N636").PasteSpecial(-4122) | Out-Null
$ws.Range("A636").Value = 45192.62309835648
$ws.Range("B636").Value = '1207dpwls@naver.com'
$ws.Range("C636").Value = '법학과'
$ws.Range("D636").Value = 20202751
$ws.Range("E636").Value = '주예진'
$ws.Range("F636").Value = '78:22'
$ws.Range("G636").Value = 0.2
$ws.Range("H636").Value = '우리나라의 국민부담률은 2010년 22.4%에서 꾸준히 상승하여 2020년 27.9%에 달하였다.'
$ws.Range("I636").Value = '38만 명'
$ws.Range("J636").Value = 0.151
$ws.Range("K636").Value = '중견기업은 신고법인수의 0.5%를 차지하는 데 부담하는 세액은 8.7%이다'
$ws.Range("L636").Value = 'Black'
$ws.Range("N636").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F636").PasteSpecial(-4122) | Out-Null

# Row 637
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A637:N637").PasteSpecial(-4122) | Out-Null
$ws.Range("A637").Value = 45192.649936828704
$ws.Range("B637").Value = 'gwkang0330@gmail.com'
$ws.Range("C637").Value = '소프트웨어학부'
$ws.Range("D637").Value = 20235101
$ws.Range("E637").Value = '강건우'
$ws.Range("F637").Value = '77:23'
$ws.Range("G637").Value = 0.15
$ws.Range("H637").Value = '조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다'
$ws.Range("I637").Value = '166만 명'
$ws.Range("J637").Value = 0.002
$ws.Range("K637").Value = '그 외 대기업은 신고법인수의 8.3%를 차지하는 데 부담하는 세액은 41.2%이다'
$ws.Range("L637").Value = 'Black'
$ws.Range("N637").Value = '국민부담률을 OECD 평균 수준으로 높여야 한다'
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F637").PasteSpecial(-4122) | Out-Null

# Row 638
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A638:N638").PasteSpecial(-4122) | Out-Null
$ws.Range("A638").Value = 45192.65257868056
$ws.Range("B638").Value = 'minjoo902@naver.com'
$ws.Range("C638").Value = '금융재무학과'
$ws.Range("D638").Value = 20192827
$ws.Range("E638").Value = '김민주'
$ws.Range("F638").Value = '74:26'
$ws.Range("G638").Value = 0.2
$ws.Range("H638").Value = '조세 뿐 아니라 사회보장기여금을 포함하는 모든 강제적 납부액을 명목 GDP 대비 비율로 표시한 것이다'
$ws.Range("I638").Value = '952만 명'
$ws.Range("J638").Value = 0.059
$ws.Range("K638").Value = '법인당 평균세액은 상호출자제한기업 > 그외 대기업 > 중견기업 > 중소기업의 순서로 많이 부담하였다.'
$ws.Range("L638").Value = 'Red'
$ws.Range("M638").Value = '국민부담률을 아일랜드 수준으로 낮춰야 한다'
$ws.Range("F3").Copy() | Out-Null
$ws.Range("F638").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("D643").Select() | Out-Null